# Actualización de datos obtenidos el 6 de abril de 2016
# - Column A ("Tipo de presupuesto") and column G ("Estado de la información")
#   are reclassified from SDMX measures to SDMX dimensions.
# - A new row (row 6) is added carrying the mapping-file references for
#   the two now-dimensional columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: concept URI prefix changes from "iaest-measure:" to "iaest-dimension:"
$ws.Range("A3").Value = "iaest-dimension:tipo-de-presupuesto"
$ws.Range("G3").Value = "iaest-dimension:estado-de-la-informacion"

# Row 4: role changes from "medida" (measure) to "dim" (dimension)
$ws.Range("A4").Value = "dim"
$ws.Range("G4").Value = "dim"

# Row 5: datatype changes from "xsd:string" to "skos:Concept"
$ws.Range("A5").Value = "skos:Concept"
$ws.Range("G5").Value = "skos:Concept"

# Row 6 (new): mapping file references for the two dimension columns
$ws.Range("A6").Value = "mapping-tipo-de-presupuesto.xlsx"
$ws.Range("G6").Value = "mapping-estado-de-la-informacion.xlsx"

# Match the formatting used throughout the rest of the table (row 5's style)
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("G5").Copy()
$ws.Range("G6").PasteSpecial(-4122)
